# Update the LR-pairs TPM-derived statistics with newly recomputed values.
# (commit: "update scripts wuth new tpm")
#
# Columns G,H,I,J  = Ligand avg/total expression + derived specificities
#                    -> function of the Sending cluster (col A)
# Columns M,N,O,P  = Receptor avg/total expression + derived specificities
#                    -> function of the Target cluster (col D)
# Columns Q,R,S,T  = Edge avg/total expression weight + derived specificities
#                    -> per-row (Sending cluster x Target cluster pair)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ligand-side values (G,H,I,J), keyed by Sending cluster (rows 2-4 ECs, 5-7 FAPs, 8-10 MuSCs) ---
$ligand = @{
    2  = @{ G = 63.91118233333333; H = 191.733547;   I = 0.4067926910433548; J = 0.4067926910433549 }
    3  = @{ G = 63.91118233333333; H = 191.733547;   I = 0.4067926910433548; J = 0.4067926910433549 }
    4  = @{ G = 63.91118233333333; H = 191.733547;   I = 0.4067926910433548; J = 0.4067926910433549 }
    5  = @{                                          I = 0.3656254573230189; J = 0.365625457323019  }
    6  = @{                                          I = 0.3656254573230189; J = 0.365625457323019  }
    7  = @{                                          I = 0.3656254573230189; J = 0.365625457323019  }
    8  = @{ G = 35.755375;         H = 107.266125;   I = 0.2275818516336261; J = 0.2275818516336262 }
    9  = @{ G = 35.755375;         H = 107.266125;   I = 0.2275818516336261; J = 0.2275818516336262 }
    10 = @{ G = 35.755375;         H = 107.266125;   I = 0.2275818516336261; J = 0.2275818516336262 }
}

# --- Receptor-side values (M,N,O,P), keyed by Target cluster (rows 2,5,8 ECs; 3,6,9 FAPs; 4,7,10 MuSCs) ---
$receptor = @{
    2  = @{ M = 568.5612486666666; N = 1705.683746;        O = 0.6737621253161296; P = 0.6737621253161296 }
    3  = @{ M = 88.00803400000001;                         O = 0.1042921587987053; P = 0.1042921587987053 }
    4  = @{ M = 187.2912243333334; N = 561.8736730000001;  O = 0.2219457158851651; P = 0.2219457158851651 }
    5  = @{ M = 568.5612486666666; N = 1705.683746;        O = 0.6737621253161296; P = 0.6737621253161296 }
    6  = @{ M = 88.00803400000001;                         O = 0.1042921587987053; P = 0.1042921587987053 }
    7  = @{ M = 187.2912243333334; N = 561.8736730000001;  O = 0.2219457158851651; P = 0.2219457158851651 }
    8  = @{ M = 568.5612486666666; N = 1705.683746;        O = 0.6737621253161296; P = 0.6737621253161296 }
    9  = @{ M = 88.00803400000001;                         O = 0.1042921587987053; P = 0.1042921587987053 }
    10 = @{ M = 187.2912243333334; N = 561.8736730000001;  O = 0.2219457158851651; P = 0.2219457158851651 }
}

# --- Edge-weight values (Q,R,S,T), one set per row ---
$edge = @{
    2  = @{ Q = 36337.421631203;  R = 327036.794680827;  S = 0.2740815080804384;  T = 0.2740815080804384  }
    3  = @{ Q = 5624.6975077722;  R = 50622.27756994979; S = 0.04242528793244624; T = 0.04242528793244624 }
    4  = @{ Q = 11970.0035878009; R = 107730.0322902081; S = 0.09028589503047017; T = 0.09028589503047017 }
    5  = @{ Q = 32660.0912316588; R = 293940.8210849292; S = 0.2463445851956391;  T = 0.2463445851956391  }
    6  = @{ Q = 5055.480700275601;                       S = 0.03813186825598155; T = 0.03813186825598156 }
    7  = @{ Q = 10758.6447158694; R = 96827.8024428246;  S = 0.08114900387139831; T = 0.08114900387139831 }
    8  = @{ Q = 20329.12065654492; R = 182962.0859089042; S = 0.153336032040052;  T = 0.153336032040052   }
    9  = @{ Q = 3146.76025868275; R = 28320.84232814475; S = 0.02373500261027753; T = 0.02373500261027754 }
    10 = @{ Q = 6696.66796024746; R = 60270.01164222713; S = 0.05051081698329658; T = 0.05051081698329658 }
}

foreach ($row in 2..10) {
    foreach ($col in $ligand[$row].Keys)   { $ws.Range("$col$row").Value = $ligand[$row][$col] }
    foreach ($col in $receptor[$row].Keys) { $ws.Range("$col$row").Value = $receptor[$row][$col] }
    foreach ($col in $edge[$row].Keys)     { $ws.Range("$col$row").Value = $edge[$row][$col] }
}
